$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 160.2659493284493
$ws.Range("C2").Value = 0.005879875738173723
$ws.Range("E2").Value = -0.05798688988061057
$ws.Range("F2").Value = -0.05999019653375162
$ws.Range("H2").Value = 154.2467948717949
